# Databases.xlsx - "Added special card ideas"
#
# 1. Cards sheet: remove the long-form "Cauterize / burn stacking" comment
#    that lived on C9 (its content now lives as a shared-string effect
#    description in the new Special Cards table), and clear that cell's
#    value.
# 2. Classes sheet: add a "Special Cards" mini-table (rows 9-15, mirrors
#    the existing Primary/STR/MAG/DEX/ENH layout from rows 1-6) plus a
#    Card/Effect lookup table in H9:I18, wrapped in a real Excel Table
#    (ListObject).
# 3. Update the active sheet / selections to match where the author was
#    last working.

$wb = $excel.ActiveWorkbook

$wsCards = $wb.Worksheets.Item("Cards")
$wsClasses = $wb.Worksheets.Item("Classes")

# --- Cards sheet: drop the burn-stacking comment on C9 and clear the cell ---
$c9 = $wsCards.Range("C9").Comment
if ($c9) {
    $c9.Delete()
}
$wsCards.Range("C9").ClearContents()

# --- Classes sheet: new "Special Cards" block (rows 9-15, cols A-E) ---
$wsClasses.Range("A9").Value2 = "Special Cards"

$wsClasses.Range("B10").Value2 = "Off"

$wsClasses.Range("A11").Value2 = "Primary"
$wsClasses.Range("B11").Value2 = "STR"
$wsClasses.Range("C11").Value2 = "MAG"
$wsClasses.Range("D11").Value2 = "DEX"
$wsClasses.Range("E11").Value2 = "ENH"

$wsClasses.Range("A12").Value2 = "STR"
$wsClasses.Range("B12").Value2 = "Stagger"
$wsClasses.Range("C12").Value2 = "Conduit"
$wsClasses.Range("D12").Value2 = "Juggernaut"
$wsClasses.Range("E12").Value2 = "Bulwark"

$wsClasses.Range("A13").Value2 = "MAG"
$wsClasses.Range("B13").Value2 = "Cauterize"
$wsClasses.Range("C13").Value2 = "Cauterize"
$wsClasses.Range("D13").Value2 = "Bound"
$wsClasses.Range("E13").Value2 = "Revive"
$wsClasses.Range("B13:E13").StyleIndex = 3

$wsClasses.Range("A14").Value2 = "DEX"
$wsClasses.Range("B14").Value2 = "Stealth"
$wsClasses.Range("C14").Value2 = "Bound"
$wsClasses.Range("D14").Value2 = "Quiver"
$wsClasses.Range("E14").Value2 = "Stealth"
$wsClasses.Range("B14:E14").StyleIndex = 3

$wsClasses.Range("A15").Value2 = "ENH"
$wsClasses.Range("B15").Value2 = "Conduit"
$wsClasses.Range("C15").Value2 = "Salvage"
$wsClasses.Range("D15").Value2 = "Revive"
$wsClasses.Range("E15").Value2 = "Stagger"
$wsClasses.Range("B15:D15").StyleIndex = 3

# --- Classes sheet: Card / Effect lookup table, H9:I18 ---
$wsClasses.Range("H9").Value2 = "Card"
$wsClasses.Range("I9").Value2 = "Effect"

$wsClasses.Range("H10").Value2 = "Bound"
$wsClasses.Range("I10").Value2 = "Leap 3 squares - Card"

$wsClasses.Range("H11").Value2 = "Bulwark"
$wsClasses.Range("I11").Value2 = "Take no damage for X turns - Card"

$wsClasses.Range("H12").Value2 = "Cauterize"
$wsClasses.Range("I12").Value2 = "Adds all burns together and deals half of total outright and applies rest as burn - Card"

$wsClasses.Range("H13").Value2 = "Juggernaut"
$wsClasses.Range("I13").Value2 = "Armor cards  50% more effective - Aura"

$wsClasses.Range("H14").Value2 = "Quiver"
$wsClasses.Range("I14").Value2 = "Can use ranged attacks twice before discarding - Aura"

$wsClasses.Range("H15").Value2 = "Revive"
$wsClasses.Range("I15").Value2 = "Heals for half of missing health - Card"

$wsClasses.Range("H16").Value2 = "Stagger"
$wsClasses.Range("I16").Value2 = "Staggers out all damage taken over 10 turns - Aura"

$wsClasses.Range("H17").Value2 = "Conduit"
$wsClasses.Range("I17").Value2 = "Plus (1% max HP) HP when discarding - Aura"

$wsClasses.Range("H18").Value2 = "Salvage"
$wsClasses.Range("I18").Value2 = "Plus (1% max AP) AP when discarding - Aura"

# Column I needs to be much wider to fit the effect text.
$wsClasses.Columns.Item(9).ColumnWidth = 76.43

# Turn H9:I18 into a real table (ListObject) with an autofilter.
$tbl = $wsClasses.ListObjects.Add(1, $wsClasses.Range("H9:I18"), $null, 1)
$tbl.Name = "Table1"

# --- View state: author ended up on Classes (was on Cards), with a
#     lingering selection of C32 left behind on Cards ---
$wsCards.Range("C32").Select()
$wsClasses.Activate()
$wsClasses.Range("I14").Select()
